$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("htru2")

$ws.Range("D8").Value = 61.2

$ws.Range("C12").Value = 91.59999999999999
$ws.Range("D12").Value = 99

$ws.Range("C13").Value = 98.3

$ws.Range("D16").Value = 97.7

$ws.Range("C17").Value = 96.59999999999999
